# Apply updates to LOT2058.xlsx as described in the diff
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LOT2058")

$ws.Range("B8").Value = '01/01/2022'
$ws.Range("C8").Value = '01/01/2022'

$ws.Range("B10").Value = 'Capacitar os alunos a utilizarem conceitos e ferramentas de matemática financeira para avaliações econômicas de projetos de engenharia'
$ws.Range("C10").Value = 'Capacitar os alunos a utilizarem conceitos e ferramentas de matemática financeira para avaliações econômicas de projetos de engenharia'

$ws.Range("B15").Value = 'Introdução à Engenharia Econômica; estimativa de custos; juros; fluxo de caixa; depreciação; comparação entre alternativa de investimentos; ponto de equilíbrio; planilha eletrônica aplicado à Engenharia Econômica.'
$ws.Range("C15").Value = 'Introdução à Engenharia Econômica; estimativa de custos; juros; fluxo de caixa; depreciação; comparação entre alternativa de investimentos; ponto de equilíbrio; planilha eletrônica aplicado à Engenharia Econômica.'

$ws.Range("B17").Value = '1.INTRODUÇÃO À ENGENHARIA ECONÔMICA: necessidades de uma análise econômica de projeto (engenharia econômica como uma ferramenta de análise de tomada de decisão); 2.ESTIMATIVA DE CUTOS: estimativa de investimento de capital (tipos de estimativas de custo de capital; índice de custos mais comuns; métodos para estimar investimento de capital); estimativa de custo de equipamentos; estimativa de custo de produção;3.JUROS: variável tempo (juros simples; juros compostos); taxa efetiva, nominal e equivalente; relações de equivalência. 4.FLUXO DE CAIXA: diagrama de fluxo de caixa; equivalência de fluxo de caixa; séries uniforme e em gradientes; elaboração de fluxo de caixa5.Depreciação: métodos de depreciação;6.COMPARAÇÃO ENTRE ALTERNATIVAS DE INVESTIMENTO: critério de Rentabilidade – método do Valor Anual Uniforme Equivalente (VAU); método do Valor Presente (VP); método da Taxa interna de Retorno (TIR); método da Taxa Interna de Retorno Modificada (TIRM); renovação e substituição de equipamentos; período de retorno;7.Ponto de Equilíbrio de um projeto (break-even-point);8.Planilhas eletrônicas aplicado à análises de Engenharia Econômica.'
$ws.Range("C17").Value = '1.INTRODUÇÃO À ENGENHARIA ECONÔMICA: necessidades de uma análise econômica de projeto (engenharia econômica como uma ferramenta de análise de tomada de decisão); 2.ESTIMATIVA DE CUTOS: estimativa de investimento de capital (tipos de estimativas de custo de capital; índice de custos mais comuns; métodos para estimar investimento de capital); estimativa de custo de equipamentos; estimativa de custo de produção;3.JUROS: variável tempo (juros simples; juros compostos); taxa efetiva, nominal e equivalente; relações de equivalência. 4.FLUXO DE CAIXA: diagrama de fluxo de caixa; equivalência de fluxo de caixa; séries uniforme e em gradientes; elaboração de fluxo de caixa5.Depreciação: métodos de depreciação;6.COMPARAÇÃO ENTRE ALTERNATIVAS DE INVESTIMENTO: critério de Rentabilidade – método do Valor Anual Uniforme Equivalente (VAU); método do Valor Presente (VP); método da Taxa interna de Retorno (TIR); método da Taxa Interna de Retorno Modificada (TIRM); renovação e substituição de equipamentos; período de retorno;7.Ponto de Equilíbrio de um projeto (break-even-point);8.Planilhas eletrônicas aplicado à análises de Engenharia Econômica.'

$ws.Range("B20").Value = 'O método utilizado tem por fundamento a Aprendizagem Baseada em Projetos (PBL) que visa desenvolver as competências técnicas relativas ao tema do projeto, bem como competências transversais, tais como: aprender a aprender, trabalho em equipe, relacionamento interpessoal, aspectos de liderança e capacidade de comunicação, dentre outras; exercícios individuais realizados no decorrer da disciplina. Para os projetos, os alunos serão divididos em grupos que desenvolverão um projeto durante o semestre relacionado a um tema de Engenharia de Econômica na área de Engenharia Bioquímica, similar ao que eles irão encontrar na vida real no efetivo exercício de sua profissão.'
$ws.Range("C20").Value = 'O método utilizado tem por fundamento a Aprendizagem Baseada em Projetos (PBL) que visa desenvolver as competências técnicas relativas ao tema do projeto, bem como competências transversais, tais como: aprender a aprender, trabalho em equipe, relacionamento interpessoal, aspectos de liderança e capacidade de comunicação, dentre outras; exercícios individuais realizados no decorrer da disciplina. Para os projetos, os alunos serão divididos em grupos que desenvolverão um projeto durante o semestre relacionado a um tema de Engenharia de Econômica na área de Engenharia Bioquímica, similar ao que eles irão encontrar na vida real no efetivo exercício de sua profissão.'

$ws.Range("B21").Value = 'A nota do projeto (NP) será individual e será a média ponderada de componentes do projeto, tais como: Projeto Preliminar, Projeto Final, envolvimento do aluno com o projeto, avaliação dos Pares, apresentação de trabalhos, dentre outros;A nota dos exercícios (NE) será a média aritmética dos exercícios individuais;NF = (0,8NP + 0,2NE)'
$ws.Range("C21").Value = 'A nota do projeto (NP) será individual e será a média ponderada de componentes do projeto, tais como: Projeto Preliminar, Projeto Final, envolvimento do aluno com o projeto, avaliação dos Pares, apresentação de trabalhos, dentre outros;A nota dos exercícios (NE) será a média aritmética dos exercícios individuais;NF = (0,8NP + 0,2NE)'

$ws.Range("B22").Value = 'A recuperação será feita por meio de uma prova escrita (PR) e a média de recuperação (MR) calculada pela fórmula: MR = (NF + PR)/2n'
$ws.Range("C22").Value = 'A recuperação será feita por meio de uma prova escrita (PR) e a média de recuperação (MR) calculada pela fórmula: MR = (NF + PR)/2n'

$ws.Range("B23").Value = '1.Blank L. e Tarquin A. Engenharia Econômica. 6ed. Mc GrawHill, 2007.2.Fraser N. M., Jewkes, Elizabeth M. Engineering Economics: Financial decision making for engineers. 5th Edition. Editora Pearson, 2013. 3.Samanez, Carlos Patrício. Engenharia Econômica. São Paulo. Editora Prentice Hall, 2009. 4.Hirschfeld, Henrique. Engenharia Econômica e análise de custos. 7ed. Editora Atalas, 2007. 5.Ehrlich, Pierre Jacques. Engenharia Econômica. São Paulo. Editora Atlas, 2005 6.Peter, MS. And Timmerhaus, K.D. Plant Design and economic for chemical engineers'
$ws.Range("C23").Value = '1.Blank L. e Tarquin A. Engenharia Econômica. 6ed. Mc GrawHill, 2007.2.Fraser N. M., Jewkes, Elizabeth M. Engineering Economics: Financial decision making for engineers. 5th Edition. Editora Pearson, 2013. 3.Samanez, Carlos Patrício. Engenharia Econômica. São Paulo. Editora Prentice Hall, 2009. 4.Hirschfeld, Henrique. Engenharia Econômica e análise de custos. 7ed. Editora Atalas, 2007. 5.Ehrlich, Pierre Jacques. Engenharia Econômica. São Paulo. Editora Atlas, 2005 6.Peter, MS. And Timmerhaus, K.D. Plant Design and economic for chemical engineers'
